$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.322531
$ws.Range("H2").Value = 0.967593
$ws.Range("I2").Value = 0.01892149513432853
$ws.Range("J2").Value = 0.01892149513432853
$ws.Range("Q2").Value = 1.084589400084667
$ws.Range("R2").Value = 9.761304600762001
$ws.Range("S2").Value = 0.01417009705524306
$ws.Range("T2").Value = 0.01417009705524305

$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.322531
$ws.Range("H3").Value = 0.967593
$ws.Range("I3").Value = 0.01892149513432853
$ws.Range("J3").Value = 0.01892149513432853
$ws.Range("O3").Value = 0.251111132886386
$ws.Range("Q3").Value = 0.3636754195873333
$ws.Range("R3").Value = 3.273078776286
$ws.Range("S3").Value = 0.004751398079085478
$ws.Range("T3").Value = 0.004751398079085476

$ws.Range("I4").Value = 0.0261208867009986
$ws.Range("J4").Value = 0.0261208867009986
$ws.Range("S4").Value = 0.01956164124951391
$ws.Range("T4").Value = 0.0195616412495139

$ws.Range("I5").Value = 0.0261208867009986
$ws.Range("J5").Value = 0.0261208867009986
$ws.Range("O5").Value = 0.251111132886386
$ws.Range("Q5").Value = 0.5020493551666667
$ws.Range("S5").Value = 0.006559245451484692
$ws.Range("T5").Value = 0.006559245451484689

$ws.Range("I6").Value = 0.954957618164673
$ws.Range("J6").Value = 0.954957618164673
$ws.Range("S6").Value = 0.7151571288088572
$ws.Range("T6").Value = 0.7151571288088571

$ws.Range("I7").Value = 0.954957618164673
$ws.Range("J7").Value = 0.954957618164673
$ws.Range("O7").Value = 0.251111132886386
$ws.Range("S7").Value = 0.2398004893558158
$ws.Range("T7").Value = 0.2398004893558158
